$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109; this shifts the existing rows
# 109..212 down to 110..213 and grows the used range to A1:R213.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new data record.
$ws.Cells.Item(109, 1).Value = 7
$ws.Cells.Item(109, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(109, 3).Value = "Ñuble"
$ws.Cells.Item(109, 4).Value = 44574
$ws.Cells.Item(109, 5).Value = 16
$ws.Cells.Item(109, 6).Value = 100112002
$ws.Cells.Item(109, 7).Value = "Pimiento"
$ws.Cells.Item(109, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 200
$ws.Cells.Item(109, 11).Value = 9500
$ws.Cells.Item(109, 12).Value = 10000
$ws.Cells.Item(109, 13).Value = 9750
$ws.Cells.Item(109, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(109, 15).Value = "Región del Maule"
$ws.Cells.Item(109, 16).Value = 650
$ws.Cells.Item(109, 17).Value = 15
$ws.Cells.Item(109, 18).Value = "Hortaliza"
